$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Total Project" (column H) values for existing rows
$hUpdates = @{
    5 = 2955
    7 = 60208
    8 = 3443
    11 = 2292
    13 = 2535
    14 = 4213
    15 = 4538
    18 = 1376
    19 = 1516
    20 = 2195
    21 = 1084
    22 = 3268
    24 = 91115
    26 = 3214
    28 = 2947
    29 = 4667
    31 = 5212
    32 = 2410
    33 = 6138
    34 = 3833
    35 = 3685
    36 = 12770
    37 = 583
    38 = 4498
    40 = 6514
    42 = 1878
    43 = 3093
    44 = 5337
    45 = 2913
    46 = 5481
    47 = 6235
    48 = 21305
    49 = 5014
    52 = 22487
    53 = 4505
    54 = 15596
    55 = 14562
    56 = 10718
    57 = 33499
    58 = 11057
    63 = 724
    64 = 155056
    65 = 9560
    66 = 2226
    68 = 8354
    69 = 8675
    71 = 3622
    78 = 7043
    79 = 3886
    80 = 3499
    82 = 12287
    84 = 16612
    86 = 2785
    87 = 6983
    91 = 1149
    92 = 1174
    94 = 978
    95 = 24069
    97 = 5642
    103 = 9023
    104 = 3121
    105 = 1085
    109 = 538
    111 = 4290
    112 = 3277
    113 = 1138
    115 = 1081
    116 = 31704
    119 = 797
    120 = 1287
    121 = 618
    122 = 7941
    123 = 2837
    124 = 4920
    125 = 421
    126 = 1274
    127 = 977
    128 = 4062
    129 = 15358
    130 = 21197
    131 = 1368
    132 = 1118
    133 = 608
    134 = 18421
    135 = 17397
    136 = 4705
    137 = 1071
    141 = 314
    143 = 931
    144 = 3891
    145 = 2423
    151 = 89712
    154 = 1329
    155 = 1143
    156 = 10104
    157 = 2097
    158 = 815
    160 = 12452
    161 = 1546
    165 = 4453
    166 = 432
    167 = 2344
    168 = 1629
    169 = 2986
    173 = 7799
    175 = 3347
    177 = 148
    178 = 1737
    182 = 2585
    187 = 12442
    188 = 1342
    189 = 31764
    191 = 5079
    195 = 1574
    197 = 963
    200 = 1906
    201 = 1141
    203 = 9808
    204 = 4975
    205 = 2052
    209 = 3705
    210 = 3655
    216 = 4205
    217 = 1466
    218 = 3316
    220 = 8658
    223 = 36993
    224 = 1670
    225 = 3942
    226 = 3903
    228 = 1468
    229 = 7963
    230 = 3083
    233 = 1103
    234 = 4051
    236 = 3598
    238 = 1668
    240 = 20870
    242 = 1357
    243 = 4793
    244 = 1582
    245 = 3404
    246 = 1339
    247 = 2589
    251 = 50432
    253 = 488
    254 = 2217
    256 = 4068
    257 = 437
    260 = 1759
    261 = 16619
    263 = 632
    264 = 7656
    266 = 1029
    267 = 10561
    268 = 3126
    271 = 4714
    272 = 4909
    274 = 11954
    275 = 7996
    277 = 1009
    280 = 1642
    282 = 24653
    284 = 4464
    288 = 1528
    289 = 696
    290 = 314
    292 = 1712
    293 = 2924
    297 = 447
    299 = 8265
    300 = 2716
    302 = 1262
    303 = 5794
    304 = 1199
    305 = 7757
    307 = 2145
    310 = 11222
    311 = 4129
    313 = 1295
    314 = 1560
    315 = 3701
    316 = 849
    317 = 816
    318 = 1087
    321 = 741
    322 = 992
    323 = 2468
    324 = 4677
    326 = 10196
    327 = 3008
    328 = 3003
    329 = 2678
    330 = 1138
    332 = 2868
    333 = 3137
    334 = 7410
    335 = 1339
    336 = 5014
    337 = 8172
    338 = 5327
    339 = 1517
    340 = 14724
    341 = 1444
    345 = 1551
    346 = 2462
    349 = 289
    351 = 22073
    356 = 13495
    357 = 2921
}

foreach ($row in $hUpdates.Keys) {
    $ws.Cells.Item($row, 8).Value = $hUpdates[$row]
}

# New row 363: liferay/liferay-portal
$ws.Cells.Item(363, 1).Value = "liferay"
$ws.Cells.Item(363, 2).Value = "liferay-portal"
$ws.Cells.Item(363, 3).Value = "Enterprise Resource Planning"
$ws.Cells.Item(363, 4).Value = "38220df7f113ab58af84b3297d81db6bb40d3320"
$ws.Cells.Item(363, 5).Value = ""
$ws.Cells.Item(363, 6).Value = "2024-10-12 00:30:11+00:00"
$ws.Cells.Item(363, 7).Value = $true
$ws.Cells.Item(363, 8).Value = 109041
$ws.Cells.Item(363, 9).Value = 6549
$ws.Cells.Item(363, 10).Value = 3519
$ws.Cells.Item(363, 11).Value = 202
$ws.Cells.Item(363, 12).Value = 1480
